# Update "想去人数" (number of people wanting to go) counts on the
# "展览" and "全部类型" sheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" -> F column row numbers and their new values
$ws1 = $wb.Sheets.Item("展览")
$updates1 = @{
    3  = 38
    4  = 572
    5  = 1793
    9  = 2156
    10 = 46
    12 = 1365
    13 = 480
    15 = 300
    17 = 9
    20 = 56
    23 = 1175
    27 = 276
    28 = 345
}
foreach ($row in $updates1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updates1[$row]
}

# Sheet "全部类型" -> F column row numbers and their new values
$ws4 = $wb.Sheets.Item("全部类型")
$updates4 = @{
    3  = 38
    4  = 572
    5  = 1793
    10 = 2156
    11 = 46
    13 = 1365
    14 = 480
    16 = 300
    18 = 9
    21 = 56
    24 = 1175
    28 = 276
    29 = 345
}
foreach ($row in $updates4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updates4[$row]
}
